# Retail & Demo comment improvements
# Adds the "Actor BitVars" sheet (bit-flag reference table for the Actor
# misc-flags byte) after "BitVars", and restores "Actors" as the active /
# selected tab (it had drifted to "Scripts" in the previous save).

$wb = $excel.ActiveWorkbook

# --- Add the new "Actor BitVars" worksheet at the end of the tab strip ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Actor BitVars"

# --- Header row (bold, matches the existing BitVars/Pins sheets' style) ---
$newSheet.Range("A1").Value = "Var"
$newSheet.Range("B1").Value = "Purpose"
$newSheet.Range("C1").Value = "Name"
$newSheet.Range("A1:C1").Font.Bold = $true

# --- Data rows: one per bit of the actor misc-flags byte ---
$data = @(
  @("0x01", "Kid is strong (Hunk-O-Matic used)", "kActorMiscFlagStrong"),
  @("0x02", "Kid is green tentacle's friend (recording contract)", "kActorMiscFlagGTFriend"),
  @("0x04", "Kid knows publisher's address (watched TV)", "kActorMiscFlagWatchedTV"),
  @("0x08", "Kid is not Weird Ed's friend", "kActorMiscFlagEdsEnemy"),
  @("0x10", "", "kActorMiscFlag_10"),
  @("0x20", "", "kActorMiscFlag_20"),
  @("0x40", "Kid stops moving", "kActorMiscFlagFreeze"),
  @("0x80", "Kid is invisible (dead or in radiation suit)", "kActorMiscFlagHide")
)

$r = 2
foreach ($row in $data) {
    $newSheet.Range("A$r").Value = $row[0]
    if ($row[1] -ne "") {
        $newSheet.Range("B$r").Value = $row[1]
    }
    $newSheet.Range("C$r").Value = $row[2]
    $r++
}

# --- Column sizing: autofit the three content columns like the rest of the workbook ---
$newSheet.Columns.Item(1).AutoFit() | Out-Null
$newSheet.Columns.Item(2).AutoFit() | Out-Null
$newSheet.Columns.Item(3).AutoFit() | Out-Null
$newSheet.Columns.Item(4).ColumnWidth = $newSheet.Columns.Item(2).ColumnWidth

# --- Page setup: portrait, matching the rest of the workbook's printed sheets ---
$newSheet.PageSetup.Orientation = 1

# --- Leave the new sheet's own selection parked like the source workbook ---
$newSheet.Range("C37").Select() | Out-Null

# --- Restore "Actors" as the active sheet / selected tab ---
$wb.Worksheets.Item("Actors").Activate()
